# "Add files via upload" — the workbook's Sheet1!A1 value was updated from 3 to 4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 4
